$wb = $excel.ActiveWorkbook

# ============================================================
# Step 1: insert the new "2022-Q4" sheet right after "总计",
# i.e. before the existing "2022-Q3" sheet.
# ============================================================
$wsTotal = $wb.Worksheets.Item(1)
$new = $wb.Worksheets.Add($null, $wsTotal)
$new.Name = "2022-Q4"

# Re-fetch sheet objects fresh (by name) now that the sheet collection
# has shifted, and use "2022-Q3" purely as a formatting template so the
# new sheet reuses the very same style indices (bold/border/centered
# header, centered index column) instead of Excel minting new-but-
# equivalent styles.
$wsQ3 = $wb.Worksheets.Item("2022-Q3")

# Header row (B1:H1) styling
$wsQ3.Range("B1:H1").Copy()
$new.Range("B1:H1").PasteSpecial(-4122)

# Index column (A2:A29) styling
$wsQ3.Range("A2:A29").Copy()
$new.Range("A2:A29").PasteSpecial(-4122)

# ============================================================
# Step 2: populate the "2022-Q4" sheet with its fund-holdings data
# ============================================================
$new.Cells.Item(1,2).Value = "基金代码"
$new.Cells.Item(1,3).Value = "基金名称"
$new.Cells.Item(1,4).Value = "基金规模"
$new.Cells.Item(1,5).Value = "股票总仓位"
$new.Cells.Item(1,6).Value = "仓位占比"
$new.Cells.Item(1,7).Value = "持有市值(亿元)"
$new.Cells.Item(1,8).Value = "仓位排名"
$new.Cells.Item(2,1).Value = 0
$new.Cells.Item(2,2).Value = "'000729"
$new.Cells.Item(2,3).Value = "建信中小盘先锋股票A"
$new.Cells.Item(2,4).Value = "'31.93"
$new.Cells.Item(2,5).Value = "'91.13"
$new.Cells.Item(2,6).Value = "'4.13"
$new.Cells.Item(2,7).Value = "'1.3187"
$new.Cells.Item(2,8).Value = 4
$new.Cells.Item(3,1).Value = 1
$new.Cells.Item(3,2).Value = "'530005"
$new.Cells.Item(3,3).Value = "建信优化配置混合A"
$new.Cells.Item(3,4).Value = "'24.05"
$new.Cells.Item(3,5).Value = "'88.67"
$new.Cells.Item(3,6).Value = "'4.10"
$new.Cells.Item(3,7).Value = "'0.9860"
$new.Cells.Item(3,8).Value = 4
$new.Cells.Item(4,1).Value = 2
$new.Cells.Item(4,2).Value = "'010003"
$new.Cells.Item(4,3).Value = "景顺长城电子信息产业股票A"
$new.Cells.Item(4,4).Value = "'18.37"
$new.Cells.Item(4,5).Value = "'93.05"
$new.Cells.Item(4,6).Value = "'4.42"
$new.Cells.Item(4,7).Value = "'0.8120"
$new.Cells.Item(4,8).Value = 7
$new.Cells.Item(5,1).Value = 3
$new.Cells.Item(5,2).Value = "'000756"
$new.Cells.Item(5,3).Value = "建信潜力新蓝筹股票A"
$new.Cells.Item(5,4).Value = "'15.30"
$new.Cells.Item(5,5).Value = "'90.58"
$new.Cells.Item(5,6).Value = "'4.02"
$new.Cells.Item(5,7).Value = "'0.6151"
$new.Cells.Item(5,8).Value = 4
$new.Cells.Item(6,1).Value = 4
$new.Cells.Item(6,2).Value = "'014967"
$new.Cells.Item(6,3).Value = "建信潜力新蓝筹股票C"
$new.Cells.Item(6,4).Value = "'10.68"
$new.Cells.Item(6,5).Value = "'90.58"
$new.Cells.Item(6,6).Value = "'4.02"
$new.Cells.Item(6,7).Value = "'0.4293"
$new.Cells.Item(6,8).Value = 4
$new.Cells.Item(7,1).Value = 5
$new.Cells.Item(7,2).Value = "'013919"
$new.Cells.Item(7,3).Value = "建信中小盘先锋股票C"
$new.Cells.Item(7,4).Value = "'10.14"
$new.Cells.Item(7,5).Value = "'91.13"
$new.Cells.Item(7,6).Value = "'4.13"
$new.Cells.Item(7,7).Value = "'0.4188"
$new.Cells.Item(7,8).Value = 4
$new.Cells.Item(8,1).Value = 6
$new.Cells.Item(8,2).Value = "'010004"
$new.Cells.Item(8,3).Value = "景顺长城电子信息产业股票C"
$new.Cells.Item(8,4).Value = "'5.90"
$new.Cells.Item(8,5).Value = "'93.05"
$new.Cells.Item(8,6).Value = "'4.42"
$new.Cells.Item(8,7).Value = "'0.2608"
$new.Cells.Item(8,8).Value = 7
$new.Cells.Item(9,1).Value = 7
$new.Cells.Item(9,2).Value = "'004476"
$new.Cells.Item(9,3).Value = "景顺长城沪港深领先科技股票"
$new.Cells.Item(9,4).Value = "'8.15"
$new.Cells.Item(9,5).Value = "'87.92"
$new.Cells.Item(9,6).Value = "'2.75"
$new.Cells.Item(9,7).Value = "'0.2241"
$new.Cells.Item(9,8).Value = 9
$new.Cells.Item(10,1).Value = 8
$new.Cells.Item(10,2).Value = "'000020"
$new.Cells.Item(10,3).Value = "景顺长城品质投资混合A"
$new.Cells.Item(10,4).Value = "'7.65"
$new.Cells.Item(10,5).Value = "'88.29"
$new.Cells.Item(10,6).Value = "'2.74"
$new.Cells.Item(10,7).Value = "'0.2096"
$new.Cells.Item(10,8).Value = 10
$new.Cells.Item(11,1).Value = 9
$new.Cells.Item(11,2).Value = "'006615"
$new.Cells.Item(11,3).Value = "工银战略新兴产业混合A"
$new.Cells.Item(11,4).Value = "'5.87"
$new.Cells.Item(11,5).Value = "'75.35"
$new.Cells.Item(11,6).Value = "'3.56"
$new.Cells.Item(11,7).Value = "'0.2090"
$new.Cells.Item(11,8).Value = 5
$new.Cells.Item(12,1).Value = 10
$new.Cells.Item(12,2).Value = "'013365"
$new.Cells.Item(12,3).Value = "汇添富产业升级混合A"
$new.Cells.Item(12,4).Value = "'4.03"
$new.Cells.Item(12,5).Value = "'85.05"
$new.Cells.Item(12,6).Value = "'3.57"
$new.Cells.Item(12,7).Value = "'0.1439"
$new.Cells.Item(12,8).Value = 8
$new.Cells.Item(13,1).Value = 11
$new.Cells.Item(13,2).Value = "'011506"
$new.Cells.Item(13,3).Value = "建信高端装备股票A"
$new.Cells.Item(13,4).Value = "'3.84"
$new.Cells.Item(13,5).Value = "'89.26"
$new.Cells.Item(13,6).Value = "'3.22"
$new.Cells.Item(13,7).Value = "'0.1236"
$new.Cells.Item(13,8).Value = 8
$new.Cells.Item(14,1).Value = 12
$new.Cells.Item(14,2).Value = "'001858"
$new.Cells.Item(14,3).Value = "建信鑫利灵活配置混合"
$new.Cells.Item(14,4).Value = "'3.07"
$new.Cells.Item(14,5).Value = "'90.23"
$new.Cells.Item(14,6).Value = "'3.61"
$new.Cells.Item(14,7).Value = "'0.1108"
$new.Cells.Item(14,8).Value = 9
$new.Cells.Item(15,1).Value = 13
$new.Cells.Item(15,2).Value = "'009598"
$new.Cells.Item(15,3).Value = "景顺长城科技创新三年定期开放灵活配置混合"
$new.Cells.Item(15,4).Value = "'2.67"
$new.Cells.Item(15,5).Value = "'91.56"
$new.Cells.Item(15,6).Value = "'3.22"
$new.Cells.Item(15,7).Value = "'0.0860"
$new.Cells.Item(15,8).Value = 8
$new.Cells.Item(16,1).Value = 14
$new.Cells.Item(16,2).Value = "'005967"
$new.Cells.Item(16,3).Value = "鹏华创新驱动混合"
$new.Cells.Item(16,4).Value = "'1.46"
$new.Cells.Item(16,5).Value = "'93.54"
$new.Cells.Item(16,6).Value = "'5.14"
$new.Cells.Item(16,7).Value = "'0.0750"
$new.Cells.Item(16,8).Value = 7
$new.Cells.Item(17,1).Value = 15
$new.Cells.Item(17,2).Value = "'000522"
$new.Cells.Item(17,3).Value = "华润元大信息传媒科技混合"
$new.Cells.Item(17,4).Value = "'1.38"
$new.Cells.Item(17,5).Value = "'62.01"
$new.Cells.Item(17,6).Value = "'4.39"
$new.Cells.Item(17,7).Value = "'0.0606"
$new.Cells.Item(17,8).Value = 6
$new.Cells.Item(18,1).Value = 16
$new.Cells.Item(18,2).Value = "'006616"
$new.Cells.Item(18,3).Value = "工银战略新兴产业混合C"
$new.Cells.Item(18,4).Value = "'1.34"
$new.Cells.Item(18,5).Value = "'75.35"
$new.Cells.Item(18,6).Value = "'3.56"
$new.Cells.Item(18,7).Value = "'0.0477"
$new.Cells.Item(18,8).Value = 5
$new.Cells.Item(19,1).Value = 17
$new.Cells.Item(19,2).Value = "'005914"
$new.Cells.Item(19,3).Value = "景顺长城智能生活混合"
$new.Cells.Item(19,4).Value = "'1.07"
$new.Cells.Item(19,5).Value = "'88.91"
$new.Cells.Item(19,6).Value = "'3.24"
$new.Cells.Item(19,7).Value = "'0.0347"
$new.Cells.Item(19,8).Value = 5
$new.Cells.Item(20,1).Value = 18
$new.Cells.Item(20,2).Value = "'011507"
$new.Cells.Item(20,3).Value = "建信高端装备股票C"
$new.Cells.Item(20,4).Value = "'0.90"
$new.Cells.Item(20,5).Value = "'89.26"
$new.Cells.Item(20,6).Value = "'3.22"
$new.Cells.Item(20,7).Value = "'0.0290"
$new.Cells.Item(20,8).Value = 8
$new.Cells.Item(21,1).Value = 19
$new.Cells.Item(21,2).Value = "'001223"
$new.Cells.Item(21,3).Value = "鹏华文化传媒娱乐股票"
$new.Cells.Item(21,4).Value = "'0.77"
$new.Cells.Item(21,5).Value = "'83.57"
$new.Cells.Item(21,6).Value = "'3.10"
$new.Cells.Item(21,7).Value = "'0.0239"
$new.Cells.Item(21,8).Value = 10
$new.Cells.Item(22,1).Value = 20
$new.Cells.Item(22,2).Value = "'001162"
$new.Cells.Item(22,3).Value = "前海开源优势蓝筹股票A"
$new.Cells.Item(22,4).Value = "'0.42"
$new.Cells.Item(22,5).Value = "'91.86"
$new.Cells.Item(22,6).Value = "'4.44"
$new.Cells.Item(22,7).Value = "'0.0186"
$new.Cells.Item(22,8).Value = 9
$new.Cells.Item(23,1).Value = 21
$new.Cells.Item(23,2).Value = "'004223"
$new.Cells.Item(23,3).Value = "金信多策略精选灵活配置混合"
$new.Cells.Item(23,4).Value = "'0.31"
$new.Cells.Item(23,5).Value = "'93.96"
$new.Cells.Item(23,6).Value = "'4.91"
$new.Cells.Item(23,7).Value = "'0.0152"
$new.Cells.Item(23,8).Value = 6
$new.Cells.Item(24,1).Value = 22
$new.Cells.Item(24,2).Value = "'013366"
$new.Cells.Item(24,3).Value = "汇添富产业升级混合C"
$new.Cells.Item(24,4).Value = "'0.28"
$new.Cells.Item(24,5).Value = "'85.05"
$new.Cells.Item(24,6).Value = "'3.57"
$new.Cells.Item(24,7).Value = "'0.0100"
$new.Cells.Item(24,8).Value = 8
$new.Cells.Item(25,1).Value = 23
$new.Cells.Item(25,2).Value = "'004931"
$new.Cells.Item(25,3).Value = "华润元大价值优选混合C"
$new.Cells.Item(25,4).Value = "'0.16"
$new.Cells.Item(25,5).Value = "'74.11"
$new.Cells.Item(25,6).Value = "'3.36"
$new.Cells.Item(25,7).Value = "'0.0054"
$new.Cells.Item(25,8).Value = 9
$new.Cells.Item(26,1).Value = 24
$new.Cells.Item(26,2).Value = "'015436"
$new.Cells.Item(26,3).Value = "建信优化配置混合C"
$new.Cells.Item(26,4).Value = "'0.11"
$new.Cells.Item(26,5).Value = "'88.67"
$new.Cells.Item(26,6).Value = "'4.10"
$new.Cells.Item(26,7).Value = "'0.0045"
$new.Cells.Item(26,8).Value = 4
$new.Cells.Item(27,1).Value = 25
$new.Cells.Item(27,2).Value = "'004930"
$new.Cells.Item(27,3).Value = "华润元大价值优选混合A"
$new.Cells.Item(27,4).Value = "'0.13"
$new.Cells.Item(27,5).Value = "'74.11"
$new.Cells.Item(27,6).Value = "'3.36"
$new.Cells.Item(27,7).Value = "'0.0044"
$new.Cells.Item(27,8).Value = 9
$new.Cells.Item(28,1).Value = 26
$new.Cells.Item(28,2).Value = "'001638"
$new.Cells.Item(28,3).Value = "前海开源优势蓝筹股票C"
$new.Cells.Item(28,4).Value = "'0.08"
$new.Cells.Item(28,5).Value = "'91.86"
$new.Cells.Item(28,6).Value = "'4.44"
$new.Cells.Item(28,7).Value = "'0.0036"
$new.Cells.Item(28,8).Value = 9
$new.Cells.Item(29,1).Value = 27
$new.Cells.Item(29,2).Value = "'016906"
$new.Cells.Item(29,3).Value = "景顺长城品质投资混合C"
$new.Cells.Item(29,4).Value = "'0.00"
$new.Cells.Item(29,5).Value = "'88.29"
$new.Cells.Item(29,6).Value = "'2.74"
$new.Cells.Item(29,7).Value = 0
$new.Cells.Item(29,8).Value = 10

# ============================================================
# Step 3: update the "总计" (summary) sheet - a new row is inserted
# for 2022-Q4 at the top of the data, pushing the existing 2022-Q3
# and 2022-Q1 rows down by one.
# ============================================================

# Row 4 is brand new -- give A4 the same index-column style as A2/A3
# before writing its value.
$wsTotal.Range("A3").Copy()
$wsTotal.Range("A4").PasteSpecial(-4122)

# Shift old row 3 (2022-Q1) down to row 4
$wsTotal.Cells.Item(4,1).Value = 2
$wsTotal.Cells.Item(4,2).Value = "2022-Q1"
$wsTotal.Cells.Item(4,3).Value = 1
$wsTotal.Cells.Item(4,4).Value = 0

# Shift old row 2 (2022-Q3) down to row 3
$wsTotal.Cells.Item(3,1).Value = 1
$wsTotal.Cells.Item(3,2).Value = "2022-Q3"
$wsTotal.Cells.Item(3,3).Value = 29
$wsTotal.Cells.Item(3,4).Value = 5.89

# New row 2: 2022-Q4 data
$wsTotal.Cells.Item(2,1).Value = 0
$wsTotal.Cells.Item(2,2).Value = "2022-Q4"
$wsTotal.Cells.Item(2,3).Value = 28
$wsTotal.Cells.Item(2,4).Value = 6.28

Write-Output "edit complete"
